{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the homework-answer-key edit described by the diff:\n//   - \"...create a Q-Q plot to be sure the distribution of sample means is\n//      normal:\" becomes \"...create a histogram to be sure the distribution\n//      of differences is normal in order to be able to claim that the\n//      distribution sample means is normal:\"\n//   - \"The data appear to be normal\" becomes\n//      \"The differences appear to be normal\"\n\nconst body = context.document.body;\n\n// --- Edit 1: \"Q-Q plot\" -> \"histogram\" -------------------------------\nconst qqHits = body.search(\"Q-Q plot\", { matchCase: true });\nqqHits.load(\"items\");\nawait context.sync();\n\nif (qqHits.items.length === 0) {\n  throw new Error('Search text \"Q-Q plot\" not found.');\n}\nqqHits.items[0].insertText(\"histogram\", \"Replace\");\n\n// --- Edit 2: \"the distribution of sample means is normal:\" -----------\n//     -> \"the distribution of differences is normal in order to be able\n//         to claim that the distribution sample means is normal:\"\nconst distHits = body.search(\"the distribution of sample means is normal:\", {\n  matchCase: true,\n});\ndistHits.load(\"items\");\nawait context.sync();\n\nif (distHits.items.length === 0) {\n  throw new Error(\n    'Search text \"the distribution of sample means is normal:\" not found.'\n  );\n}\ndistHits.items[0].insertText(\n  \"the distribution of differences is normal in order to be able to claim that the distribution sample means is normal:\",\n  \"Replace\"\n);\n\n// --- Edit 3: \"The data appear to be normal\" -> \"The differences appear\n//     to be normal\" -------------------------------------------------\nconst dataHits = body.search(\"The data appear to be normal\", {\n  matchCase: true,\n});\ndataHits.load(\"items\");\nawait context.sync();\n\nif (dataHits.items.length === 0) {\n  throw new Error('Search text \"The data appear to be normal\" not found.');\n}\ndataHits.items[0].insertText(\"The differences appear to be normal\", \"Replace\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n#\n# Applies the homework-answer-key edit described by the diff:\n#   - \"...create a Q-Q plot to be sure the distribution of sample means is\n#      normal:\" becomes \"...create a histogram to be sure the distribution\n#      of differences is normal in order to be able to claim that the\n#      distribution sample means is normal:\"\n#   - \"The data appear to be normal\" becomes\n#      \"The differences appear to be normal\"\n\n$d = $word.ActiveDocument\n\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nfunction Replace-DocText($FindText, $ReplaceWith) {\n    $rng = $d.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $rng.Find.Text = $FindText\n    $rng.Find.Replacement.Text = $ReplaceWith\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = $wdFindContinue\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n\n    $ok = $rng.Find.Execute($FindText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $ReplaceWith, $wdReplaceAll)\n    if (-not $ok) {\n        throw \"Find/Replace failed: '$FindText' was not found.\"\n    }\n}\n\n# 1. \"Q-Q plot\" -> \"histogram\"\nReplace-DocText \"Q-Q plot\" \"histogram\"\n\n# 2. \"the distribution of sample means is normal:\" -> the longer sentence\nReplace-DocText \"the distribution of sample means is normal:\" \"the distribution of differences is normal in order to be able to claim that the distribution sample means is normal:\"\n\n# 3. \"The data appear to be normal\" -> \"The differences appear to be normal\"\nReplace-DocText \"The data appear to be normal\" \"The differences appear to be normal\"\n\n$d.Save()\n"}
